# data set balancing - ACER
# Extends the labelled-comment dataset on Sheet1: fills in the "text"/"flagged"
# columns for a few already-present rows (118-120), then appends 20 brand-new
# rows (121-140, ids COM120..COM139) continuing the same B/C/D/E/F/G layout,
# five of which (121-124) also get text + flagged-word annotations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: row 118 - add the "text" comment (F column). A/B/C/D/E already hold
# data (id COM117, all-zero labels) and are left untouched.
# ---------------------------------------------------------------------------
$ws.Range("F117").Copy()
$ws.Range("F118").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F118").Value = "ammige hinawama thiyena lassana puthek"
$ws.Rows.Item(118).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 2: append the 20 new ids (COM120..COM139) into A121:A140, with the
# B/C/D/E label columns defaulted to 0 (most stay that way; a few below get
# updated to reflect flagged content).
# ---------------------------------------------------------------------------
for ($row = 121; $row -le 140; $row++) {
    $num = 120 + ($row - 121)
    $ws.Cells.Item($row, 1).Value = "COM$num"
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = 0
}

# ---------------------------------------------------------------------------
# Step 3: row 119 - add the "text" comment for COM118.
# ---------------------------------------------------------------------------
$ws.Range("F117").Copy()
$ws.Range("F119").PasteSpecial(-4122)
$ws.Range("F119").Value = "Budu saranai puthe oyata"
$ws.Rows.Item(119).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 4: row 120 (COM119) - mark as offensive/hate-flagged, add text +
# flagged word.
# ---------------------------------------------------------------------------
$ws.Cells.Item(120, 2).Value = 1       # hate
$ws.Cells.Item(120, 3).Value = 1       # offensive
$ws.Cells.Item(120, 4).Value = 0       # neither
$ws.Cells.Item(120, 5).Value = 1       # class

$ws.Range("F117").Copy()
$ws.Range("F120").PasteSpecial(-4122)
$ws.Range("F120").Value = "Carima panditayek wei mu issarahata"
$ws.Range("G120").Value = "carima"
$ws.Rows.Item(120).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 5: row 121 (COM120) - offensive, add text + flagged word.
# ---------------------------------------------------------------------------
$ws.Cells.Item(121, 2).Value = 0
$ws.Cells.Item(121, 3).Value = 1
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 1

$ws.Range("F117").Copy()
$ws.Range("F121").PasteSpecial(-4122)
$ws.Range("F121").Value = "Ane me pakayage nambar Eka nadda"
$ws.Range("G121").Value = "pakayage"
$ws.Rows.Item(121).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 6: row 122 (COM121) - not flagged, but text + flagged word recorded.
# ---------------------------------------------------------------------------
$ws.Cells.Item(122, 2).Value = 0
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 0

$ws.Range("F117").Copy()
$ws.Range("F122").PasteSpecial(-4122)
$ws.Range("F122").Value = "Pissu pakayek muta riport gahapalla okkom ekathu.wela"
$ws.Range("G122").Value = "pakayek"
$ws.Rows.Item(122).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 7: row 123 (COM122) - not flagged, just a text comment (no flagged
# word captured).
# ---------------------------------------------------------------------------
$ws.Cells.Item(123, 2).Value = 0
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 0

$ws.Range("F117").Copy()
$ws.Range("F123").PasteSpecial(-4122)
$ws.Range("F123").Value = "Waddek"
$ws.Rows.Item(123).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 8: row 124 (COM123) - offensive, add text + flagged word.
# ---------------------------------------------------------------------------
$ws.Cells.Item(124, 2).Value = 0
$ws.Cells.Item(124, 3).Value = 1
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 1

$ws.Range("F117").Copy()
$ws.Range("F124").PasteSpecial(-4122)
$ws.Range("F124").Value = "Ithin laccna pacyo dn kt whn plyn"
$ws.Range("G124").Value = "pacyo"
$ws.Rows.Item(124).RowHeight = 16.5

# ---------------------------------------------------------------------------
# Step 9: restore the view - scroll position + active cell, matching where
# the author was working when they saved.
# ---------------------------------------------------------------------------
try {
    $excel.ActiveWindow.ScrollRow = 109
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
$ws.Range("F125").Select()

Write-Output "edit applied"
